# Auto-generated edit script: update Universalis market snapshot values
# per the commit diff (per-cell static value updates; workbook has no formulas).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 2719.2
$ws.Range("I107").Value = 2770.4285
$ws.Range("J107").Value = 2599.6667
$ws.Range("K107").Value = 2770.4285
$ws.Range("L107").Value = 2599.6667
$ws.Range("M107").Value = -850.4285
$ws.Range("N107").Value = -6439.6667
$ws.Range("H111").Value = 2066.6843
$ws.Range("J111").Value = 1455.8
$ws.Range("L111").Value = 4367.4
$ws.Range("N111").Value = -10501.4
$ws.Range("H137").Value = 2042373
$ws.Range("J137").Value = 4547633.5
$ws.Range("L137").Value = 13642900.5
$ws.Range("N137").Value = -13648000.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 870.53125
$ws.Range("I2").Value = 695.4231
$ws.Range("K2").Value = 695.4231
$ws.Range("M2").Value = -582.4231
$ws.Range("H32").Value = 142925.31
$ws.Range("I32").Value = 164227.1
$ws.Range("K32").Value = 164227.1
$ws.Range("M32").Value = -163940.1
$ws.Range("H45").Value = 103808.7
$ws.Range("I45").Value = 127386
$ws.Range("J45").Value = 9499.5
$ws.Range("K45").Value = 127386
$ws.Range("L45").Value = 9499.5
$ws.Range("M45").Value = -127009
$ws.Range("N45").Value = -10253.5
$ws.Range("H63").Value = 12077.625
$ws.Range("I63").Value = 3227.6
$ws.Range("J63").Value = 26827.666
$ws.Range("K63").Value = 3227.6
$ws.Range("L63").Value = 26827.666
$ws.Range("M63").Value = -2541.6
$ws.Range("N63").Value = -28199.666
$ws.Range("H66").Value = 12077.625
$ws.Range("I66").Value = 3227.6
$ws.Range("J66").Value = 26827.666
$ws.Range("K66").Value = 16138
$ws.Range("L66").Value = 134138.33
$ws.Range("M66").Value = -12706
$ws.Range("N66").Value = -141002.33
$ws.Range("H116").Value = 870.53125
$ws.Range("I116").Value = 695.4231
$ws.Range("K116").Value = 695.4231
$ws.Range("M116").Value = 1598.5769
$ws.Range("H122").Value = 2697.9524
$ws.Range("I122").Value = 2317.8333
$ws.Range("K122").Value = 6953.499899999999
$ws.Range("M122").Value = -4503.499899999999
$ws.Range("H132").Value = 7198.769
$ws.Range("I132").Value = 7264.8335
$ws.Range("K132").Value = 21794.5005
$ws.Range("M132").Value = -19264.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 870.53125
$ws.Range("I3").Value = 695.4231
$ws.Range("K3").Value = 695.4231
$ws.Range("M3").Value = -581.4231
$ws.Range("H94").Value = 1141.1777
$ws.Range("I94").Value = 1086.6177
$ws.Range("K94").Value = 1086.6177
$ws.Range("M94").Value = -635.6177
$ws.Range("H134").Value = 20002894
$ws.Range("I134").Value = 2451.5925
$ws.Range("K134").Value = 7354.7775
$ws.Range("M134").Value = -4819.7775

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4180
$ws.Range("I31").Value = 1891.7273
$ws.Range("K31").Value = 1891.7273
$ws.Range("M31").Value = -1596.7273
$ws.Range("H34").Value = 4180
$ws.Range("I34").Value = 1891.7273
$ws.Range("K34").Value = 1891.7273
$ws.Range("M34").Value = -1689.7273
$ws.Range("H132").Value = 9010875
$ws.Range("I132").Value = 1571.75
$ws.Range("J132").Value = 25643436
$ws.Range("K132").Value = 4715.25
$ws.Range("L132").Value = 76930308
$ws.Range("M132").Value = -2185.25
$ws.Range("N132").Value = -76935368

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 344.72726
$ws.Range("I12").Value = 123.333336
$ws.Range("K12").Value = 370.000008
$ws.Range("M12").Value = -197.000008
$ws.Range("H74").Value = 18381.8
$ws.Range("I74").Value = 19945
$ws.Range("J74").Value = 17991
$ws.Range("K74").Value = 59835
$ws.Range("L74").Value = 53973
$ws.Range("N74").Value = -56095
$ws.Range("M74").Value = -58774
$ws.Range("H77").Value = 18381.8
$ws.Range("I77").Value = 19945
$ws.Range("J77").Value = 17991
$ws.Range("K77").Value = 179505
$ws.Range("L77").Value = 161919
$ws.Range("N77").Value = -172527
$ws.Range("M77").Value = -174201
$ws.Range("H107").Value = 166667020
$ws.Range("J107").Value = 333333500
$ws.Range("L107").Value = 1000000500
$ws.Range("N107").Value = -1000004340
$ws.Range("H117").Value = 19612980
$ws.Range("J117").Value = 30309204
$ws.Range("L117").Value = 90927612
$ws.Range("N117").Value = -90934496
$ws.Range("H129").Value = 17208866
$ws.Range("I129").Value = 2443.3845
$ws.Range("J129").Value = 29635726
$ws.Range("K129").Value = 7330.1535
$ws.Range("L129").Value = 88907178
$ws.Range("M129").Value = -2330.1535
$ws.Range("N129").Value = -88917178
$ws.Range("H131").Value = 6771029
$ws.Range("I131").Value = 6995970.5
$ws.Range("J131").Value = 6599015.5
$ws.Range("K131").Value = 20987911.5
$ws.Range("L131").Value = 19797046.5
$ws.Range("M131").Value = -20982871.5
$ws.Range("N131").Value = -19807126.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 54999.5
$ws.Range("H70").Value = 7714.656
$ws.Range("I70").Value = 8311.574000000001
$ws.Range("J70").Value = 5710.7144
$ws.Range("K70").Value = 8311.574000000001
$ws.Range("L70").Value = 5710.7144
$ws.Range("M70").Value = -8041.574000000001
$ws.Range("N70").Value = -6250.7144
$ws.Range("H72").Value = 54999.5
$ws.Range("H73").Value = 7714.656
$ws.Range("I73").Value = 8311.574000000001
$ws.Range("J73").Value = 5710.7144
$ws.Range("K73").Value = 8311.574000000001
$ws.Range("L73").Value = 5710.7144
$ws.Range("M73").Value = -7375.574000000001
$ws.Range("N73").Value = -7582.7144
$ws.Range("H80").Value = 15681189
$ws.Range("I80").Value = 66244
$ws.Range("J80").Value = 50034068
$ws.Range("K80").Value = 66244
$ws.Range("L80").Value = 50034068
$ws.Range("M80").Value = -65246
$ws.Range("N80").Value = -50036064
$ws.Range("H83").Value = 15681189
$ws.Range("I83").Value = 66244
$ws.Range("J83").Value = 50034068
$ws.Range("K83").Value = 331220
$ws.Range("L83").Value = 250170340
$ws.Range("M83").Value = -326228
$ws.Range("N83").Value = -250180324
$ws.Range("H126").Value = 2753.3333
$ws.Range("I126").Value = 2753.3333
$ws.Range("K126").Value = 8259.999899999999
$ws.Range("M126").Value = -5789.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 44114
$ws.Range("J63").Value = 46333.332
$ws.Range("L63").Value = 46333.332
$ws.Range("N63").Value = -47831.332
$ws.Range("H66").Value = 44114
$ws.Range("J66").Value = 46333.332
$ws.Range("L66").Value = 138999.996
$ws.Range("N66").Value = -146487.996
$ws.Range("H130").Value = 57131
$ws.Range("J130").Value = 57131
$ws.Range("L130").Value = 57131
$ws.Range("N130").Value = -67171
$ws.Range("H132").Value = 9725.467000000001
$ws.Range("I132").Value = 2764.889
$ws.Range("J132").Value = 20166.334
$ws.Range("K132").Value = 8294.667000000001
$ws.Range("L132").Value = 60499.00199999999
$ws.Range("M132").Value = -5764.667000000001
$ws.Range("N132").Value = -65559.00199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 744.14813
$ws.Range("I113").Value = 880.44446
$ws.Range("K113").Value = 2641.33338
$ws.Range("M113").Value = -471.33338
$ws.Range("H122").Value = 1682.3055
$ws.Range("I122").Value = 1763.5385
$ws.Range("J122").Value = 1471.1
$ws.Range("K122").Value = 5290.6155
$ws.Range("L122").Value = 4413.299999999999
$ws.Range("M122").Value = -2840.6155
$ws.Range("N122").Value = -9313.299999999999
